$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (old C/D/E shift right to D/E/F)
$ws.Columns("C:C").Insert()

# New "Id" header in C1, matching the existing header formatting (bold white on grey fill)
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Id"
$ws.Range("C1").NumberFormat = "@"

# New signup id values beneath the header (text-formatted column, mixed numeric/text content)
$ws.Range("C2").Value = 1
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C3").Value = "x"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C4").Value = 3
$ws.Range("C4").NumberFormat = "@"

# Give the new column roughly the same width as the other hidden id column (B)
$ws.Columns("C:C").ColumnWidth = 12.44140625

# Rebuild the autofilter over the new A1:F1 range
$ws.AutoFilterMode = $false
$ws.Range("A1:F1").AutoFilter()

# Keep the _FilterDatabase defined name in sync with the autofilter range
foreach ($n in $wb.Names) {
    $n.RefersTo = "=Sheet1!`$A`$1:`$F`$1"
}

# Move the active selection the way the author left it
$ws.Range("C2").Select()
